$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: header row for the new "RAGear" block (mirrors row 3's formatting) ---
$ws.Range("A3:K3").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null

$ws.Range("A8").Value = "RAGear"
$ws.Range("B8").Value = "pulley"
$ws.Range("C8").Value = "ratio"
$ws.Range("D8").Value = 32
$ws.Range("E8").Value = 63.683950000000003
$ws.Range("F8").Value = "FULL STEP"
$ws.Range("G8").Value = "HALF STEP"
$ws.Range("H8").Value = "min/step"
$ws.Range("I8").Value = "FULL STEP"
$ws.Range("J8").Value = "HALF STEP"
$ws.Range("K8").Value = "min/step"
# D8 keeps the default (no special border) style, unlike D3
$ws.Range("D8").Borders.Item(7).LineStyle = 0
$ws.Range("D8").Borders.Item(8).LineStyle = 0
$ws.Range("D8").Borders.Item(9).LineStyle = 0
$ws.Range("D8").Borders.Item(10).LineStyle = 0

# --- Row 9: first ratio row for RAGear (mirrors row 6's formatting/formulas, thick bottom) ---
$ws.Range("E6:K6").Copy() | Out-Null
$ws.Range("E9").PasteSpecial(-4122) | Out-Null

$ws.Range("A9").Value = 180
$ws.Range("B9").Value = 20
$ws.Range("C9").Formula = "=B9/A9"
$ws.Range("E9").Formula = "=E`$4"
$ws.Range("F9").Formula = "=E9/360/C9"
$ws.Range("G9").Formula = "=1/F9"
$ws.Range("H9").Formula = "=G9*60"
$ws.Range("I9").Formula = "=F9*2"
$ws.Range("J9").Formula = "=1/I9"
$ws.Range("K9").Formula = "=J9*60"

# --- Row 10: second ratio row for RAGear (mirrors row 6's formatting/formulas, thick bottom) ---
$ws.Range("E6:K6").Copy() | Out-Null
$ws.Range("E10").PasteSpecial(-4122) | Out-Null

$ws.Range("A10").Value = 180
$ws.Range("B10").Value = 16
$ws.Range("C10").Formula = "=B10/A10"
$ws.Range("E10").Formula = "=E`$4"
$ws.Range("F10").Formula = "=E10/360/C10"
$ws.Range("G10").Formula = "=1/F10"
$ws.Range("H10").Formula = "=G10*60"
$ws.Range("I10").Formula = "=F10*2"
$ws.Range("J10").Formula = "=1/I10"
$ws.Range("K10").Formula = "=J10*60"

# --- Update selection to match the new active cell ---
$excel.CutCopyMode = $false
$ws.Range("B11").Select() | Out-Null

$wb.Save()
